$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Detach the old "summary" row (old row 5) and "footer" row (old row 6)
#    from their merges so the rows can be rebuilt as normal data rows; their
#    content will be re-created later at rows 18 and 19. Clearing the old
#    anchor cells (only possible once un-merged) drops the now-unused shared
#    strings so that, when the same text is written back near the end of the
#    script, it is appended as a fresh shared-string entry instead of
#    re-using the old slot.
# ---------------------------------------------------------------------------
$ws.Range("K5:N5").UnMerge()
$ws.Range("A6:E6").UnMerge()
$ws.Range("F6:G6").UnMerge()
$ws.Range("I6:N6").UnMerge()
$ws.Range("K5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("I6").ClearContents()

# ---------------------------------------------------------------------------
# 2. Fill in row 4 (already existing, previously blank) with the first
#    transaction line.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS."
$ws.Range("H4").Value = "0:0"
$ws.Range("L4").Value = 114
$ws.Range("N4").Value = 1

# ---------------------------------------------------------------------------
# 3. Propagate row 4's formatting down to the 13 new data rows (5-17) in one
#    shot so every column keeps the same font/border/fill as the template
#    row, then populate each one with its own data.
# ---------------------------------------------------------------------------
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
    @{R=5;  N=2;  Name="AUGMENTIN 457MG/5ML SUSP. 70 ML";          Ratio="1:0";    Amt=137;    Pct=1;    H=25.5}
    @{R=6;  N=3;  Name="BLOKATENS 10/160MG 28 F.C.TABS.";          Ratio="0:0";    Amt=160;    Pct=1;    H=24.75}
    @{R=7;  N=4;  Name="COLOVATIL 30 F.C. TABS";                   Ratio="0:0";    Amt=63;     Pct=1;    H=25.5}
    @{R=8;  N=5;  Name="GAVISCON LIQUID 24 SACHETS 10 ML";         Ratio="0:9";    Amt=12;     Pct=0.04; H=25.5}
    @{R=9;  N=6;  Name="GINKGO BILOBA 30 CAPS.";                   Ratio="0:0";    Amt=186;    Pct=1;    H=24.75}
    @{R=10; N=7;  Name="MILGA ADVANCE 30 F.C. TABS";                Ratio="0:0";    Amt=136.5;  Pct=1;    H=25.5}
    @{R=11; N=8;  Name="PERLOC 40MG 14 F.C.TAB.";                   Ratio="0:0";    Amt=68.25;  Pct=1;    H=24.75}
    @{R=12; N=9;  Name="RHINEX 0.05% INFANTILE NASAL DROPS 10 ML";  Ratio="2:0";    Amt=18;     Pct=1;    H=25.5}
    @{R=13; N=10; Name="RIVO 320MG 20*10 TABS";                     Ratio="1:2";    Amt=14.1;   Pct=0.1;  H=25.5}
    @{R=14; N=11; Name="VASTAREL MR 35MG 30 F.C.TAB.";              Ratio="2:0";    Amt=175;    Pct=1;    H=24.75}
    @{R=15; N=12; Name="WATER FOR INJECTION AMP. 5 ML";             Ratio="7816:0"; Amt=2.5;    Pct=1;    H=25.5}
    @{R=16; N=13; Name="سويت كوكو";                                  Ratio="22:0";   Amt=25;     Pct=1;    H=24.75}
    @{R=17; N=14; Name="مرطب شفاه لونا جوز هند ابيض";                Ratio="3:0";    Amt=20;     Pct=1;    H=25.5}
)

$ws.Rows.Item(4).RowHeight = 24.75

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.N
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 8).Value = $row.Ratio
    $ws.Cells.Item($r, 12).Value = $row.Amt
    $ws.Cells.Item($r, 14).Value = $row.Pct
    $ws.Rows.Item($r).RowHeight = $row.H
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}

# ---------------------------------------------------------------------------
# 4. Re-create the totals row (previously row 5) now at row 18.
# ---------------------------------------------------------------------------
$ws.Range("K5:N5").Copy()
$ws.Range("K18:N18").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K18").Value = 1131.3499999999999
$ws.Rows.Item(18).RowHeight = 25.5
$ws.Range("K18:N18").Merge()

# ---------------------------------------------------------------------------
# 5. Re-create the footer row (previously row 6) now at row 19.
# ---------------------------------------------------------------------------
$ws.Range("A6:E6").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F6:G6").Copy()
$ws.Range("F19:G19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H6").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I6:N6").Copy()
$ws.Range("I19:N19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A19").Value = "Monday, 5 January, 2026 10:09 AM"
$ws.Range("F19").Value = "1/1"
$ws.Range("I19").Value = "developed by : Abdelaziz Talaat"
$ws.Rows.Item(19).RowHeight = 17.25

$ws.Range("A19:E19").Merge()
$ws.Range("F19:G19").Merge()
$ws.Range("I19:N19").Merge()

# ---------------------------------------------------------------------------
# 6. Clear the now-stale content that used to live directly in rows 5 and 6
#    (it has been re-created above at rows 18/19); only needed where the new
#    data rows (5-17) do not already overwrite those cells.
# ---------------------------------------------------------------------------
